# 22/09/2017 HARISH CHICK IN
#
# 1) Merge the two runs that together spell out
#    "Thu Sep 19 11:11:52 PDT 2017" into a single run.
# 2) Append a new "Fri Sep 21 11:42:38 PDT 2017" purchase-details entry
#    right after the "Thu Sep 19" entry's "Amount balance - 5800.0" line.

$d = $word.ActiveDocument

# --- Change 1: merge the split date/time runs -----------------------------
$null = $d.Content.Find.Execute(
    "Thu Sep 19 11:11:52 PDT 2017", $false, $false, $false, $false, $false,
    $true, 1, $false, "Thu Sep 19 11:11:52 PDT 2017", 2)

# --- Change 2: insert the new "Fri Sep 21" block ---------------------------
# Locate the paragraph that ends the "Thu Sep 19" entry
# (the bold "Amount balance ... - 5800.0" line).
$anchorIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*Amount balance*- 5800.0*") {
        $anchorIndex = $i
    }
}

$tab = [char]9

# Content (and bold state) for each new paragraph, in order.
$newParas = @(
    @{ Bold = $true;  Text = $null },
    @{ Bold = $false; Text = "Fri Sep 21 11:42:38 PDT 2017" },
    @{ Bold = $false; Text = "Person Name" + $tab + $tab + $tab + $tab + "- ST" },
    @{ Bold = $false; Text = "---------------------------------------------------------------" },
    @{ Bold = $false; Text = "Item Name" + $tab + $tab + $tab + $tab + "- CARROT EVE" },
    @{ Bold = $false; Text = "Number of Pockets" + $tab + $tab + $tab + "- 1" },
    @{ Bold = $false; Text = "Number of KGs" + $tab + $tab + $tab + "- 90" },
    @{ Bold = $false; Text = "Rate" + $tab + $tab + $tab + $tab + $tab + "- 22" },
    @{ Bold = $false; Text = "Total Price" + $tab + $tab + $tab + $tab + "- 1980.0" },
    @{ Bold = $true;  Text = "Amount balance" + $tab + $tab + $tab + "- 7780.0" },
    @{ Bold = $false; Text = $null },
    @{ Bold = $true;  Text = $null }
)

# Insert all the (initially empty) paragraphs first, chaining off the
# previous one so they land in the right order right after $anchorIndex.
$anchor = $d.Paragraphs.Item($anchorIndex).Range
for ($i = 0; $i -lt $newParas.Count; $i++) {
    $null = $anchor.InsertParagraphAfter()
    $anchor = $d.Paragraphs.Item($anchorIndex + 1 + $i).Range
}

# Now fill in text / formatting for each newly created paragraph.
for ($i = 0; $i -lt $newParas.Count; $i++) {
    $p = $d.Paragraphs.Item($anchorIndex + 1 + $i)
    $spec = $newParas[$i]

    $p.Range.Font.Bold = $spec.Bold
    if ($spec.Text -ne $null) {
        $p.Range.Text = $spec.Text
        $p.Range.Font.Bold = $spec.Bold
    }
}
